# Reverses the order of the "Periodo Mora" (col E) and "Valor Mora" (col F)
# values for the worker detail rows (16..82) on Hoja1 - i.e. flips the block
# top-to-bottom so the periods run in ascending order (1808 .. 2402) instead
# of descending (2402 .. 1808), carrying their paired "Valor Mora" amount
# along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow = 82

$eVals = @()
$fVals = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $eVals += , $ws.Cells.Item($r, 5).Value2
    $fVals += , $ws.Cells.Item($r, 6).Value2
}

$count = $eVals.Count
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $srcIdx = $count - 1 - $i
    $ws.Cells.Item($r, 5).Value2 = $eVals[$srcIdx]
    $ws.Cells.Item($r, 6).Value2 = $fVals[$srcIdx]
}
